# Auto-generated: updates currentAveragePrice / Leve profit columns (H-N)
# across all 8 sheets, refreshed from the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 795
$ws.Range("J12").Value = 790
$ws.Range("L12").Value = 790
$ws.Range("N12").Value = -1130
$ws.Range("H19").Value = 2516.4
$ws.Range("I19").Value = 2371.7856
$ws.Range("K19").Value = 2371.7856
$ws.Range("M19").Value = -2196.7856
$ws.Range("H40").Value = 13448.65
$ws.Range("I40").Value = 4328.8335
$ws.Range("J40").Value = 17357.143
$ws.Range("K40").Value = 4328.8335
$ws.Range("L40").Value = 17357.143
$ws.Range("M40").Value = -4153.8335
$ws.Range("N40").Value = -17707.143
$ws.Range("H113").Value = 2379.625
$ws.Range("I113").Value = 2463
$ws.Range("J113").Value = 1796
$ws.Range("K113").Value = 2463
$ws.Range("L113").Value = 1796
$ws.Range("M113").Value = 791
$ws.Range("N113").Value = -8304
$ws.Range("H137").Value = 3756.4443
$ws.Range("I137").Value = 3024.9412
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 9074.8236
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = -6524.8236
$ws.Range("N137").Value = -20100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7597.2915
$ws.Range("I61").Value = 5639.357
$ws.Range("K61").Value = 5639.357
$ws.Range("M61").Value = -5427.357
$ws.Range("H74").Value = 2517.4614
$ws.Range("I74").Value = 2090.3
$ws.Range("K74").Value = 2090.3
$ws.Range("M74").Value = -1216.3
$ws.Range("H77").Value = 2517.4614
$ws.Range("I77").Value = 2090.3
$ws.Range("K77").Value = 10451.5
$ws.Range("M77").Value = -6083.5
$ws.Range("H88").Value = 2198.5
$ws.Range("J88").Value = 2924.5
$ws.Range("L88").Value = 2924.5
$ws.Range("N88").Value = -3736.5
$ws.Range("H91").Value = 2198.5
$ws.Range("J91").Value = 2924.5
$ws.Range("L91").Value = 2924.5
$ws.Range("N91").Value = -5732.5
$ws.Range("H97").Value = 2671.7646
$ws.Range("I97").Value = 1067.6666
$ws.Range("J97").Value = 4476.375
$ws.Range("K97").Value = 1067.6666
$ws.Range("L97").Value = 4476.375
$ws.Range("M97").Value = -571.6666
$ws.Range("N97").Value = -5468.375
$ws.Range("H136").Value = 7597.2915
$ws.Range("I136").Value = 5639.357
$ws.Range("K136").Value = 16918.071
$ws.Range("M136").Value = -14368.071

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2164.375
$ws.Range("I134").Value = 1996.8636
$ws.Range("K134").Value = 5990.5908
$ws.Range("M134").Value = -3455.5908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 930.0909
$ws.Range("I19").Value = 2037.6666
$ws.Range("J19").Value = 514.75
$ws.Range("K19").Value = 2037.6666
$ws.Range("L19").Value = 514.75
$ws.Range("M19").Value = -1867.6666
$ws.Range("N19").Value = -854.75
$ws.Range("H24").Value = 930.0909
$ws.Range("I24").Value = 2037.6666
$ws.Range("J24").Value = 514.75
$ws.Range("K24").Value = 2037.6666
$ws.Range("L24").Value = 514.75
$ws.Range("M24").Value = -1867.6666
$ws.Range("N24").Value = -854.75
$ws.Range("H31").Value = 6967.517
$ws.Range("I31").Value = 7164.4165
$ws.Range("J31").Value = 6022.4
$ws.Range("K31").Value = 7164.4165
$ws.Range("L31").Value = 6022.4
$ws.Range("M31").Value = -6869.4165
$ws.Range("N31").Value = -6612.4
$ws.Range("H34").Value = 6967.517
$ws.Range("I34").Value = 7164.4165
$ws.Range("J34").Value = 6022.4
$ws.Range("K34").Value = 7164.4165
$ws.Range("L34").Value = 6022.4
$ws.Range("M34").Value = -6962.4165
$ws.Range("N34").Value = -6426.4
$ws.Range("H86").Value = 8179.3335
$ws.Range("I86").Value = 6975.6665
$ws.Range("K86").Value = 6975.6665
$ws.Range("M86").Value = -5852.6665
$ws.Range("H89").Value = 8179.3335
$ws.Range("I89").Value = 6975.6665
$ws.Range("K89").Value = 34878.3325
$ws.Range("M89").Value = -29262.3325
$ws.Range("H96").Value = 26555.7
$ws.Range("J96").Value = 26555.7
$ws.Range("L96").Value = 26555.7
$ws.Range("N96").Value = -32047.7
$ws.Range("H105").Value = 2827.2
$ws.Range("I105").Value = 2568
$ws.Range("K105").Value = 2568
$ws.Range("M105").Value = -821
$ws.Range("H130").Value = 94998
$ws.Range("J130").Value = 94998
$ws.Range("L130").Value = 94998
$ws.Range("N130").Value = -105038
$ws.Range("H141").Value = 86339.60000000001
$ws.Range("J141").Value = 86339.60000000001
$ws.Range("L141").Value = 86339.60000000001
$ws.Range("N141").Value = -96699.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 200007310
$ws.Range("I8").Value = 200007310
$ws.Range("K8").Value = 600021930
$ws.Range("M8").Value = -600021791
$ws.Range("H18").Value = 306
$ws.Range("I18").Value = 306
$ws.Range("K18").Value = 918
$ws.Range("M18").Value = -749
$ws.Range("H44").Value = 583.875
$ws.Range("I44").Value = 656.2
$ws.Range("K44").Value = 1968.6
$ws.Range("M44").Value = -1570.6
$ws.Range("H50").Value = 612.1177
$ws.Range("I50").Value = 348.9
$ws.Range("J50").Value = 988.1429000000001
$ws.Range("K50").Value = 1046.7
$ws.Range("L50").Value = 2964.4287
$ws.Range("M50").Value = -565.6999999999998
$ws.Range("N50").Value = -3926.4287
$ws.Range("H53").Value = 612.1177
$ws.Range("I53").Value = 348.9
$ws.Range("J53").Value = 988.1429000000001
$ws.Range("K53").Value = 1046.7
$ws.Range("L53").Value = 2964.4287
$ws.Range("M53").Value = -565.6999999999998
$ws.Range("N53").Value = -3926.4287
$ws.Range("H92").Value = 434.14285
$ws.Range("J92").Value = 450
$ws.Range("L92").Value = 1350
$ws.Range("N92").Value = -3846
$ws.Range("H117").Value = 2028
$ws.Range("I117").Value = 289
$ws.Range("J117").Value = 2607.6667
$ws.Range("K117").Value = 867
$ws.Range("L117").Value = 7823.000100000001
$ws.Range("M117").Value = 2575
$ws.Range("N117").Value = -14707.0001
$ws.Range("H134").Value = 4476
$ws.Range("I134").Value = 597.5
$ws.Range("J134").Value = 19990
$ws.Range("K134").Value = 1792.5
$ws.Range("L134").Value = 59970
$ws.Range("M134").Value = 3277.5
$ws.Range("N134").Value = -70110

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 20001.25
$ws.Range("I20").Value = 5005
$ws.Range("K20").Value = 5005
$ws.Range("M20").Value = -4760
$ws.Range("H70").Value = 3527.4443
$ws.Range("I70").Value = 3502.7222
$ws.Range("J70").Value = 3576.889
$ws.Range("K70").Value = 3502.7222
$ws.Range("L70").Value = 3576.889
$ws.Range("M70").Value = -3232.7222
$ws.Range("N70").Value = -4116.889
$ws.Range("H73").Value = 3527.4443
$ws.Range("I73").Value = 3502.7222
$ws.Range("J73").Value = 3576.889
$ws.Range("K73").Value = 3502.7222
$ws.Range("L73").Value = 3576.889
$ws.Range("M73").Value = -2566.7222
$ws.Range("N73").Value = -5448.889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3408.3333
$ws.Range("I22").Value = 3112.5
$ws.Range("K22").Value = 3112.5
$ws.Range("M22").Value = -2817.5
$ws.Range("H27").Value = 3408.3333
$ws.Range("I27").Value = 3112.5
$ws.Range("K27").Value = 3112.5
$ws.Range("M27").Value = -3005.5
$ws.Range("H33").Value = 29000
$ws.Range("I33").Value = 28000
$ws.Range("K33").Value = 28000
$ws.Range("M33").Value = -27710
$ws.Range("H55").Value = 771.125
$ws.Range("I55").Value = 466.1
$ws.Range("K55").Value = 466.1
$ws.Range("M55").Value = -293.1
$ws.Range("H132").Value = 45102.47
$ws.Range("I132").Value = 47733.875
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 143201.625
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -140671.625
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 3930.25
$ws.Range("I136").Value = 1648
$ws.Range("J136").Value = 5299.6
$ws.Range("K136").Value = 4944
$ws.Range("L136").Value = 15898.8
$ws.Range("M136").Value = -2394
$ws.Range("N136").Value = -20998.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 18131.889
$ws.Range("I45").Value = 16800
$ws.Range("J45").Value = 18298.375
$ws.Range("K45").Value = 16800
$ws.Range("L45").Value = 18298.375
$ws.Range("M45").Value = -16309
$ws.Range("N45").Value = -19280.375
